$d = $word.ActiveDocument

# Word's internal "line break" character (matches a <w:br/> when scanning Find text
# across run boundaries; using it lets a single Find/Replace span - and merge - runs).
$br = [char]11

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Medical report - doctor phone
Replace-Text "Lékař: [[PERSON_1]], [[PHONE_1]], e-mail: [[EMAIL_1]]" "Lékař: [[PERSON_1]], tel.: [[AMOUNT_1]], e-mail: [[EMAIL_1]]"

# 2. Bank report - department phone
Replace-Text "Oddělení: [[PERSON_3]]: [[PERSON_4]], e-mail: [[EMAIL_2]], [[PHONE_2]]" "Oddělení: [[PERSON_3]]: [[PERSON_4]], e-mail: [[EMAIL_2]], tel.: [[AMOUNT_2]]"

# 3. Bank report - limit amount
Replace-Text "Kontrola účtů č. [[BANK_1]], [[BANK_2]] a [[BANK_3]] prokázala několik odchozích plateb přesahujících interní limit 500 000 Kč. Upozorněno na nesrovnalosti u klienta: [[PERSON_5]], bytem: Antonínská 8, Brno. Doporučeno předat věc k dalšímu šetření oddělení AML ([[PERSON_6]])." "Kontrola účtů č. [[BANK_1]], [[BANK_2]] a [[BANK_3]] prokázala několik odchozích plateb přesahujících interní limit [[AMOUNT_3]]. Upozorněno na nesrovnalosti u klienta: [[PERSON_5]], bytem: Antonínská 8, Brno. Doporučeno předat věc k dalšímu šetření oddělení AML ([[PERSON_6]])."

# 4. Invoice - accounting advisory amounts
Replace-Text "Účetní poradenství — 20 hodin × 1 000 Kč/h = 20 000 Kč" "Účetní poradenství — 20 hodin × [[AMOUNT_4]]/h = [[AMOUNT_5]]"

# 5. Invoice - travel costs (structural: split into two w:t runs joined by a w:br)
Replace-Text "Cestovní náklady — 2 000 [[PERSON_8]] k úhradě: 22 000 Kč" ("Cestovní náklady — [[AMOUNT_6]]" + $br + "Celkem k úhradě: [[AMOUNT_7]]")

# 5b. Invoice - processed by (renumbered PERSON_9 -> PERSON_8)
Replace-Text "Zpracovala: [[PERSON_9]], fakturační oddělení" "Zpracovala: [[PERSON_8]], fakturační oddělení"

# 6. HR record - employee renumbered PERSON_10 -> PERSON_9
Replace-Text "Společnost: AUTOCOMP Global a.s., [[ICO_5]], [[ADDRESS_6]]: [[PERSON_10]], nar. 22.07.1988, RČ: [[BIRTH_ID_2]], pozice: Projektový manažer" "Společnost: AUTOCOMP Global a.s., [[ICO_5]], [[ADDRESS_6]]: [[PERSON_9]], nar. 22.07.1988, RČ: [[BIRTH_ID_2]], pozice: Projektový manažer"

# 7. HR record - evaluation PERSON_11 -> PERSON_10
Replace-Text "[[PERSON_11]] dosáhl vynikajících výsledků v projektu E-Drive 2.0, zvýšení efektivity o 12 %." "[[PERSON_10]] dosáhl vynikajících výsledků v projektu E-Drive 2.0, zvýšení efektivity o 12 %."

# 8. HR record - bonus amount
Replace-Text "Doporučeno k ročnímu bonusu 25 000 Kč." "Doporučeno k ročnímu bonusu [[AMOUNT_8]]."

# 9. HR record - dept head PERSON_12 -> PERSON_11
Replace-Text "Vedoucí oddělení: [[PERSON_12]], podpis: M. Konečný." "Vedoucí oddělení: [[PERSON_11]], podpis: M. Konečný."

# 10. IT security - author PERSON_13 -> PERSON_12, phone -> amount
Replace-Text "Autor: [[PERSON_13]], CISO, e-mail: [[EMAIL_4]], [[PHONE_3]]" "Autor: [[PERSON_12]], CISO, e-mail: [[EMAIL_4]], tel.: +420 [[AMOUNT_9]]"

# 11. IT security - IP address
Replace-Text "Událost: Detekován neautorizovaný přístup do interní sítě z IP 185.63.115.42 dne 18. 8. 2025 v 02:37. Incident se týká uživatelského účtu: [[EMAIL_5]] (oddělení vývoje)." "Událost: Detekován neautorizovaný přístup do interní sítě z IP [[IP_1]] dne 18. 8. 2025 v 02:37. Incident se týká uživatelského účtu: [[EMAIL_5]] (oddělení vývoje)."

# 12. IT security - mode/lead renumbered PERSON_14 -> PERSON_13, PERSON_15 -> PERSON_14
Replace-Text "Dočasně zablokován přístup, aktivován režim [[PERSON_14]]. O případu informován vedoucí vývoje [[PERSON_15]]." "Dočasně zablokován přístup, aktivován režim [[PERSON_13]]. O případu informován vedoucí vývoje [[PERSON_14]]."

# 13. Tax report - processed by renumbered PERSON_16 -> PERSON_15
Replace-Text "Zpracoval: [[PERSON_16]], daňová poradkyně č. 3324, e-mail: [[EMAIL_6]]" "Zpracoval: [[PERSON_15]], daňová poradkyně č. 3324, e-mail: [[EMAIL_6]]"

# 14. Tax report - client renumbered PERSON_17 -> PERSON_16, phone -> amount
Replace-Text "Klient: [[PERSON_17]], OSVČ, [[ICO_8]], [[ADDRESS_7]], [[PHONE_4]]" "Klient: [[PERSON_16]], OSVČ, [[ICO_8]], [[ADDRESS_7]], tel.: [[AMOUNT_10]]"

# 15. Tax report - merge two runs ("Dle dohody..." run + separate <w:br/>"Základ daně..." run)
#     into a single run, with amounts tagged. Matching across the line-break (the $br char)
#     makes Find/Replace span both runs, collapsing them into one.
$oldCombined = "Dle dohody z 14. 2. 2025 žádá klient o odklad daňového přiznání do 1. 7. 2025." + $br + "Základ daně: 920 000 Kč. Doporučeno zaplatit zálohu ve výši 60 000 Kč do 15. 6. 2025."
$newCombined = "Dle dohody z 14. 2. 2025 žádá klient o odklad daňového přiznání do 1. 7. 2025." + $br + "Základ daně: [[AMOUNT_11]]. Doporučeno zaplatit zálohu ve výši [[AMOUNT_12]] do 15. 6. 2025."
Replace-Text $oldCombined $newCombined

# 16. Tax report - signature renumbered PERSON_16 -> PERSON_15
Replace-Text "Podepsala: Bc. [[PERSON_16]]" "Podepsala: Bc. [[PERSON_15]]"

# 17. Board meeting - participants renumbered
Replace-Text "Účastníci: [[PERSON_18]] (CEO), [[PERSON_19]] (CFO), [[PERSON_20]] (HR), [[PERSON_21]] (CTO)" "Účastníci: [[PERSON_17]] (CEO), [[PERSON_18]] (CFO), [[PERSON_19]] (HR), [[PERSON_20]] (CTO)"

# 18. Board meeting - Ostrava lead renumbered PERSON_22 -> PERSON_21
Replace-Text "Restrukturalizace pobočky Ostrava (vedoucí [[PERSON_22]])" "Restrukturalizace pobočky Ostrava (vedoucí [[PERSON_21]])"

# 19. Board meeting - minutes taker renumbered PERSON_23 -> PERSON_22
Replace-Text "Zapsala: [[PERSON_23]], asistentka představenstva" "Zapsala: [[PERSON_22]], asistentka představenstva"

# 20. Email - renumbered PERSON_24 -> PERSON_23
Replace-Text "Potřebuji, aby [[PERSON_24]] doplnil finální data do tabulky do středy." "Potřebuji, aby [[PERSON_23]] doplnil finální data do tabulky do středy."

# 21. Email signature renumbered PERSON_25 -> PERSON_24
Replace-Text ("Díky," + $br + "[[PERSON_25]]" + $br + "Projektová manažerka") ("Díky," + $br + "[[PERSON_24]]" + $br + "Projektová manažerka")

# 22. Payment - bank address renumbered PERSON_26 -> PERSON_25
Replace-Text "Banka: Komerční banka a.s., [[ICO_10]], [[PERSON_26]] 33, Praha 1" "Banka: Komerční banka a.s., [[ICO_10]], [[PERSON_25]] 33, Praha 1"

# 23. Payment - sender renumbered PERSON_27 -> PERSON_26
Replace-Text "Odesílatel: [[PERSON_27]], bytem: [[ADDRESS_9]]" "Odesílatel: [[PERSON_26]], bytem: [[ADDRESS_9]]"

# 24. Payment - amount tagged
Replace-Text "Částka: 18 250 Kč" "[[AMOUNT_13]] Kč"

# 25. Insurance - client renumbered PERSON_28 -> PERSON_27, PERSON_29 -> PERSON_28
Replace-Text "Klient: [[PERSON_28]], nar. 17.11.1985, RČ: [[BIRTH_ID_3]], bytem: Sokolská 14, [[PERSON_29]]: Životní pojištění „Bez starostí“" "Klient: [[PERSON_27]], nar. 17.11.1985, RČ: [[BIRTH_ID_3]], bytem: Sokolská 14, [[PERSON_28]]: Životní pojištění „Bez starostí“"

# 26. Insurance - insured amount tagged, PERSON_30 -> PERSON_29
Replace-Text "Pojištěná částka: 2 000 000 [[PERSON_30]] pojištění: 1. 8. 2025" "Pojištěná [[AMOUNT_14]] [[PERSON_29]] pojištění: 1. 8. 2025"

# 27. Insurance - signature renumbered PERSON_28 -> PERSON_27
Replace-Text "Podepsala: Mgr. [[PERSON_28]], podpis ručně" "Podepsala: Mgr. [[PERSON_27]], podpis ručně"
